$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 261
$wsExhibition.Range("F4").Value = 9723
$wsExhibition.Range("F5").Value = 662
$wsExhibition.Range("F6").Value = 172
$wsExhibition.Range("F8").Value = 374
$wsExhibition.Range("F9").Value = 428
$wsExhibition.Range("F13").Value = 12347
$wsExhibition.Range("F19").Value = 248
$wsExhibition.Range("F23").Value = 164
$wsExhibition.Range("F29").Value = 2155
$wsExhibition.Range("F30").Value = 1041
$wsExhibition.Range("F31").Value = 4217
$wsExhibition.Range("F32").Value = 3700
$wsExhibition.Range("F33").Value = 657
$wsExhibition.Range("F37").Value = 1336
$wsExhibition.Range("F39").Value = 777
$wsExhibition.Range("F40").Value = 33
$wsExhibition.Range("F42").Value = 448
$wsExhibition.Range("F43").Value = 574

# Sheet "演出" (Performances) - update counts and mark one event as sold out
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F12").Value = 27
$wsPerformance.Range("F16").Value = 11
$wsPerformance.Range("F18").Value = 7
$wsPerformance.Range("G3").Value = "不可售"

# Sheet "全部类型" (All types) - update "想去人数" (want-to-go count)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 261
$wsAll.Range("F5").Value = 9723
$wsAll.Range("F6").Value = 662
$wsAll.Range("F8").Value = 172
$wsAll.Range("F10").Value = 374
$wsAll.Range("F11").Value = 428
$wsAll.Range("F15").Value = 12347
$wsAll.Range("F19").Value = 248
$wsAll.Range("F23").Value = 164
$wsAll.Range("F28").Value = 2155
$wsAll.Range("F29").Value = 1041
$wsAll.Range("F30").Value = 4217
$wsAll.Range("F31").Value = 3700
$wsAll.Range("F32").Value = 657
$wsAll.Range("F36").Value = 1336
$wsAll.Range("F38").Value = 777
$wsAll.Range("F39").Value = 33
$wsAll.Range("F41").Value = 448
$wsAll.Range("F43").Value = 574

